# Auto-generated edit script applying numeric updates to Leve profit tables
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 10403
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 10403
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 10403
$ws.Range("N13").Value = -10741
# Row 129
$ws.Range("H129").Value = 2454.7163
$ws.Range("I129").Value = 6385.5293
$ws.Range("J129").Value = 1118.24
$ws.Range("K129").Value = 19156.5879
$ws.Range("L129").Value = 3354.72
$ws.Range("M129").Value = -14156.5879
$ws.Range("N129").Value = -13354.72
# Row 132
$ws.Range("H132").Value = 6255702
$ws.Range("I132").Value = 6584902.5
$ws.Range("J132").Value = 888
$ws.Range("K132").Value = 19754707.5
$ws.Range("L132").Value = 2664
$ws.Range("M132").Value = -19752177.5
$ws.Range("N132").Value = -7724

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1955.08
$ws.Range("I32").Value = 1881.0537
$ws.Range("J32").Value = 2938.5715
$ws.Range("K32").Value = 1881.0537
$ws.Range("L32").Value = 2938.5715
$ws.Range("M32").Value = -1594.0537
$ws.Range("N32").Value = -3512.5715
# Row 95
$ws.Range("H95").Value = 32000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 32000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 32000
$ws.Range("N95").Value = -37492
# Row 97
$ws.Range("H97").Value = 39746.848
$ws.Range("I97").Value = 56525
$ws.Range("J97").Value = 1996
$ws.Range("K97").Value = 56525
$ws.Range("L97").Value = 1996
$ws.Range("M97").Value = -56029
$ws.Range("N97").Value = -2988
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 101
$ws.Range("H101").Value = 35000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 35000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490
# Row 105
$ws.Range("H105").Value = 76150.74000000001
$ws.Range("I105").Value = 47724.953
$ws.Range("J105").Value = 201224.2
$ws.Range("K105").Value = 47724.953
$ws.Range("L105").Value = 201224.2
$ws.Range("M105").Value = -45977.953
$ws.Range("N105").Value = -204718.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 38162.977
$ws.Range("I31").Value = 2491.4
$ws.Range("J31").Value = 50053.5
$ws.Range("K31").Value = 2491.4
$ws.Range("L31").Value = 50053.5
$ws.Range("M31").Value = -2196.4
$ws.Range("N31").Value = -50643.5
# Row 34
$ws.Range("H34").Value = 38162.977
$ws.Range("I34").Value = 2491.4
$ws.Range("J34").Value = 50053.5
$ws.Range("K34").Value = 2491.4
$ws.Range("L34").Value = 50053.5
$ws.Range("M34").Value = -2289.4
$ws.Range("N34").Value = -50457.5
# Row 35
$ws.Range("H35").Value = 8500
$ws.Range("I35").Value = 5000
$ws.Range("J35").Value = 12000
$ws.Range("K35").Value = 5000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = -4706
# Row 74
$ws.Range("H74").Value = 20559.25
$ws.Range("I74").Value = 3750
$ws.Range("J74").Value = 26162.334
$ws.Range("K74").Value = 3750
$ws.Range("L74").Value = 26162.334
$ws.Range("M74").Value = -2876
$ws.Range("N74").Value = -27910.334
# Row 77
$ws.Range("H77").Value = 20559.25
$ws.Range("I77").Value = 3750
$ws.Range("J77").Value = 26162.334
$ws.Range("K77").Value = 11250
$ws.Range("L77").Value = 78487.00199999999
$ws.Range("M77").Value = -6882
$ws.Range("N77").Value = -87223.00199999999
# Row 81
$ws.Range("H81").Value = 36101.8
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 36101.8
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 36101.8
$ws.Range("N81").Value = -38097.8
# Row 84
$ws.Range("H84").Value = 36101.8
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 36101.8
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 108305.4
$ws.Range("N84").Value = -118289.4
# Row 132
$ws.Range("H132").Value = 2728.06
$ws.Range("I132").Value = 2535.7856
$ws.Range("J132").Value = 3737.5
$ws.Range("K132").Value = 7607.3568
$ws.Range("L132").Value = 11212.5
$ws.Range("M132").Value = -5077.3568
$ws.Range("N132").Value = -16272.5
# Row 134
$ws.Range("H134").Value = 992.60974
$ws.Range("I134").Value = 967.1081
$ws.Range("J134").Value = 1228.5
$ws.Range("K134").Value = 2901.3243
$ws.Range("L134").Value = 3685.5
$ws.Range("M134").Value = -366.3243000000002
$ws.Range("N134").Value = -8755.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 1890
$ws.Range("I19").Value = 1890
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 5670
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -5496
# Row 23
$ws.Range("H23").Value = 544.05884
$ws.Range("I23").Value = 290
$ws.Range("J23").Value = 622.2308
$ws.Range("K23").Value = 870
$ws.Range("L23").Value = 1866.6924
$ws.Range("M23").Value = -635
$ws.Range("N23").Value = -2336.6924
# Row 74
$ws.Range("H74").Value = 5442.857
$ws.Range("I74").Value = 2833.3333
$ws.Range("J74").Value = 7400
$ws.Range("K74").Value = 8499.999899999999
$ws.Range("L74").Value = 22200
$ws.Range("M74").Value = -7438.999899999999
$ws.Range("N74").Value = -24322
# Row 77
$ws.Range("H77").Value = 5442.857
$ws.Range("I77").Value = 2833.3333
$ws.Range("J77").Value = 7400
$ws.Range("K77").Value = 25499.9997
$ws.Range("L77").Value = 66600
$ws.Range("M77").Value = -20195.9997
$ws.Range("N77").Value = -77208
# Row 121
$ws.Range("H121").Value = 1275.5
$ws.Range("I121").Value = 1599.4
$ws.Range("J121").Value = 951.6
$ws.Range("K121").Value = 4798.200000000001
$ws.Range("L121").Value = 2854.8
$ws.Range("M121").Value = -3488.200000000001
$ws.Range("N121").Value = -5474.8
# Row 131
$ws.Range("H131").Value = 772
$ws.Range("I131").Value = 359.85715
$ws.Range("J131").Value = 839.093
$ws.Range("K131").Value = 1079.57145
$ws.Range("L131").Value = 2517.279
$ws.Range("M131").Value = 3960.42855
$ws.Range("N131").Value = -12597.279
# Row 132
$ws.Range("H132").Value = 2766.75
$ws.Range("I132").Value = 1302
$ws.Range("J132").Value = 3255
$ws.Range("K132").Value = 11718
$ws.Range("L132").Value = 29295
$ws.Range("M132").Value = -9188
$ws.Range("N132").Value = -34355
# Row 138
$ws.Range("H138").Value = 7906.647
$ws.Range("I138").Value = 9416.462
$ws.Range("J138").Value = 2999.75
$ws.Range("K138").Value = 28249.386
$ws.Range("L138").Value = 8999.25
$ws.Range("M138").Value = -23109.386
$ws.Range("N138").Value = -19279.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 32
$ws.Range("H32").Value = 10150
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 19800
$ws.Range("K32").Value = 500
$ws.Range("L32").Value = 19800
$ws.Range("M32").Value = -183
$ws.Range("N32").Value = -20434
# Row 132
$ws.Range("H132").Value = 2349.6086
$ws.Range("I132").Value = 2431.2683
$ws.Range("J132").Value = 1680
$ws.Range("K132").Value = 7293.804900000001
$ws.Range("L132").Value = 5040
$ws.Range("M132").Value = -4763.804900000001
$ws.Range("N132").Value = -10100
# Row 136
$ws.Range("H136").Value = 2633.3333
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2633.3333
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 7899.999899999999
$ws.Range("M136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 150
$ws.Range("I10").Value = 150
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 150
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 19
# Row 136
$ws.Range("H136").Value = 447.66037
$ws.Range("I136").Value = 304.52
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 913.5599999999999
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = 1636.44
$ws.Range("N136").Value = -13599.9999
# Row 137
$ws.Range("H137").Value = 45143
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 45143
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 45143
$ws.Range("N137").Value = -55343

Write-Output "Applied all Leve profit updates."